$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/value assignments (safe: Excel keeps these as Text) ---
$ws.Range("D2").Value = '67.878.44'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '3.788.40'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("D7").Value = '3.785.55'
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E11").Value = '  +8.75%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").Value = '4.422.12'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '3.776.71'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '67.832.34'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("E21").Value = '  -0.91%  '
$ws.Range("E22").Value = '  -5.00%  '
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("E25").Value = '  -4.52%  '
$ws.Range("E26").Value = '  -1.99%  '
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").Value = '3.935.96'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("E32").Value = '  -7.46%  '
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("E38").Value = '  +4.97%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("E41").Value = '  -5.17%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("E45").Value = '  -2.03%  '
$ws.Range("E46").Value = '  +2.75%  '
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("E51").Value = '  -6.02%  '

# --- Numeric-looking price values that must remain Text cells ---
# Use a text formula then Copy + PasteSpecial(Values) to collapse it to a
# static value while keeping the cell type as Text (avoids Excel
# auto-converting a plain "0.996"-style literal into a Number).
$ws.Range("D4").Formula = '="0.996"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("D5").Formula = '="601.69"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="162.80"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="6.85"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="0.445"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D14").Formula = '="34.99"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="18.16"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="7.00"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="458.80"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D22").Formula = '="9.42"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="83.17"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D29").Formula = '="9.91"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D34").Formula = '="28.98"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D35").Formula = '="0.999"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="8.90"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="0.145"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D40").Formula = '="0.978"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="3.19"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="0.999"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="43.81"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="47.14"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="152.28"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D47").Formula = '="0.293"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="1.37"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="8.27"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="26.63"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$excel.CutCopyMode = $false
